$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks so relationship IDs can be rebuilt cleanly for the new row layout
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-10-18 12:33:02'
$ws.Range("B2").Value = '【業務自動化×補助金対応】生成AI活用/日本人モデル画像生成歓迎'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5405834'
$ws.Range("G2").Value = 395
$ws.Range("H2").Value = '🔥AI,Ai ◆自動化'

# Row 3
$ws.Range("A3").Value = '2025-10-18 12:33:02'
$ws.Range("B3").Value = '注目 AIプロンプトエンジニア/応答生成トレーナー募集(モバイルアプリ向け)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5415842'
$ws.Range("G3").Value = 330
$ws.Range("H3").Value = '🔥AI,Ai ◇アプリ'

# Row 4
$ws.Range("A4").Value = '2025-10-18 12:33:02'
$ws.Range("B4").Value = '【急募】ebayAPIを活用したShippingポリシー設定の専門家募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5415908'
$ws.Range("G4").Value = 183
$ws.Range("H4").Value = '🔥API'

# Row 5
$ws.Range("A5").Value = '2025-10-18 12:33:02'
$ws.Range("B5").Value = '仮想通貨取引のBOT作成'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5415610'
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = '★bot'

# Row 6
$ws.Range("A6").Value = '2025-10-18 12:33:02'
$ws.Range("B6").Value = '【メンタルヘルス】支援アプリ開発パートナー募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5415859'
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = '◆開発 ◇アプリ'

# Row 7
$ws.Range("A7").Value = '2025-10-18 12:33:02'
$ws.Range("B7").Value = '【クリエイティブ】Aurora Creative Lab 外注パートナー募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5415615'
$ws.Range("G7").Value = 18

# Row 8
$ws.Range("A8").Value = '2025-10-18 12:33:02'
$ws.Range("B8").Value = '限定公開 限定公開の仕事'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5415804'
$ws.Range("G8").Value = 13

# Row 9
$ws.Range("A9").Value = '2025-10-18 12:33:02'
$ws.Range("B9").Value = '【急募】独自ドメインのメール送信エラー解消をお願いいたします'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '~ 5,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5415841'
$ws.Range("G9").Value = 10

# Recreate hyperlinks for URL column (F) in row order, and normalize style to the Hyperlink cell style
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5405834') | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5415842') | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5415908') | Out-Null
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5415610') | Out-Null
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5415859') | Out-Null
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5415615') | Out-Null
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5415804') | Out-Null
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5415841') | Out-Null
$ws.Range("F9").Style = "Hyperlink"
